# "added new way to carry out exams"
# - Remove the "Matric" column (C): Name/Email columns remain.
# - Add a new student row: aa / aa@a.com (with mailto hyperlink, same
#   "Hyperlink" style as the existing Email cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the whole "Matric" column (shifts nothing else since it's the last col).
$ws.Columns("C").Delete()

# New row of data.
$ws.Range("A4").Value = "aa"
$ws.Range("B4").Value = "aa@a.com"

# Hyperlink the new email cell, matching the style used by the other emails.
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:aa@a.com")
$ws.Range("B4").Style = "Hyperlink"

# Resize column B to fit its (now longer) contents, and leave the new row
# selected, matching the saved state of the workbook.
$ws.Columns("B").AutoFit()
$ws.Range("B4").Select()
